$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right before the "总计" sheet
#    (so it becomes sheet #6, and "总计" shifts to #7).
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Copy header/style formatting from an existing, identically-shaped
# sheet ("2021-Q4") so the new sheet matches the workbook's look
# (bold header row + bordered index column) without inventing new
# style entries.
$srcSheet = $wb.Worksheets.Item("2021-Q4")

$srcSheet.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$srcSheet.Range("A2:A3").Copy()
$q1.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# The B:G columns hold text-like values (codes/names/decimal strings
# stored as text, same convention used by every other sheet in this
# workbook), so force Text format before writing them to avoid
# numeric coercion / loss of leading zeros, then clear the formatting
# back off afterwards (the source data cells carry no explicit style -
# only the header row and the index column A do).
$q1.Range("B2:G3").NumberFormat = "@"

# Row 2
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "519918"
$q1.Range("C2").Value = "华夏兴和混合"
$q1.Range("D2").Value = "53.07"
$q1.Range("E2").Value = "82.18"
$q1.Range("F2").Value = "3.12"
$q1.Range("G2").Value = "1.6558"
$q1.Range("H2").Value = 10

# Row 3
$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "004641"
$q1.Range("C3").Value = "万家量化睿选灵活配置混合"
$q1.Range("D3").Value = "0.16"
$q1.Range("E3").Value = "85.90"
$q1.Range("F3").Value = "1.36"
$q1.Range("G3").Value = "0.0022"
$q1.Range("H3").Value = 6

$q1.Range("B2:G3").ClearFormats()

# ------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: add a new top data row for
#    2022-Q1 (pushing the existing quarters down by one row) and
#    renumber the leading index column (A) to stay 0..n.
#    Shift the existing 5 rows down manually (bottom-up) using
#    Value2 (plain read/write, no formatting side effects) rather
#    than Rows.Insert (which drags stray formatting into the freed
#    row and bloats the style table).
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

for ($r = 6; $r -ge 2; $r--) {
    $nr = $r + 1
    $bVal = $total.Range("B$r").Value2
    $cVal = $total.Range("C$r").Value2
    $dVal = $total.Range("D$r").Value2
    $total.Range("A$nr").Value2 = ($r - 1)
    $total.Range("B$nr").Value2 = $bVal
    $total.Range("C$nr").Value2 = $cVal
    $total.Range("D$nr").Value2 = $dVal
}

# Row 7 is a brand-new cell range (the sheet used to stop at row 6) -
# clone the index-column style (bold + border) from the row above.
$total.Range("A6").Copy()
$total.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Range("A2").Value2 = 0
$total.Range("B2").Value2 = "2022-Q1"
$total.Range("C2").Value2 = 2
$total.Range("D2").Value2 = 1.66
$total.Range("A7").Value2 = 5
